$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(7,2,5,0),
    @(4,0,6,2),
    @(4,3,4,0),
    @(5,0,2,3),
    @(3,1,6,2),
    @(5,2,5,0),
    @(6,3,4,0),
    @(6,0,6,2),
    @(3,0,4,3),
    @(3,2,4,0),
    @(4,2,4,0),
    @(3,0,2,2),
    @(5,2,5,0),
    @(6,2,6,0),
    @(5,2,5,1),
    @(3,0,3,3),
    @(5,0,3,2),
    @(3,0,3,3),
    @(6,0,7,2),
    @(5,3,2,0),
    @(5,2,5,1),
    @(5,0,3,2),
    @(4,0,3,2),
    @(4,0,6,2)
)

$startRow = 1409
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# Move the active selection to the next empty row, same as Excel would
# leave it after the last data entry (matches the author's saved view).
$nextRow = $startRow + $data.Count
$ws.Range("A" + $nextRow).Select()

